# Realestate Update resale numbers 2024-01-08 18:25
# Appends a new data row (row 35) to the CityResaleNum sheet with the
# resale-number snapshot captured on 2024-01-08 at 18:25:40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 35

# Columns A and D hold values that *look* numeric/date-like ("2024-01-08",
# "01") but must be stored as literal text, matching every other row in
# this sheet. Assigning such strings straight to Range.Value lets Excel's
# smart-entry parser reinterpret them (date serials / numbers, dropping
# the leading zero), so instead we build them as text formulas in a
# scratch cell and Copy/PasteSpecial the result onto the target cell -
# exactly like typing ="2024-01-08" and converting it to its value - which
# keeps the cell a plain text cell with no special number formatting.
$scratch = $ws.Cells.Item(1000, 26)

function Set-TextValue($cell, [string]$text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy($cell)
}

Set-TextValue $ws.Cells.Item($newRow, 1) "2024-01-08"
$ws.Cells.Item($newRow, 2).Value = "18:25:40"
$ws.Cells.Item($newRow, 3).Value = "Monday"
Set-TextValue $ws.Cells.Item($newRow, 4) "01"

$scratch.Clear()

$ws.Cells.Item($newRow, 5).Value = 139531
$ws.Cells.Item($newRow, 6).Value = 142970
$ws.Cells.Item($newRow, 7).Value = 172453
$ws.Cells.Item($newRow, 8).Value = 147255
$ws.Cells.Item($newRow, 9).Value = -1
$ws.Cells.Item($newRow, 10).Value = 118220
$ws.Cells.Item($newRow, 11).Value = 224756
$ws.Cells.Item($newRow, 12).Value = 249807
$ws.Cells.Item($newRow, 13).Value = 185100
$ws.Cells.Item($newRow, 14).Value = 110388
$ws.Cells.Item($newRow, 15).Value = 40654
$ws.Cells.Item($newRow, 16).Value = 30805
$ws.Cells.Item($newRow, 17).Value = 72434
$ws.Cells.Item($newRow, 18).Value = -1
$ws.Cells.Item($newRow, 19).Value = 42141
$ws.Cells.Item($newRow, 20).Value = -1
